# Refresh the cryptos price/volume table (columns D "Price" and E
# "Volume(1h)") for rows 2-51, matching the GitHub Actions data pull.
#
# Note: several Price values look like plain decimals (e.g. "61.08",
# "1.00") which Excel's COM layer would otherwise silently coerce to a
# number when assigned via .Value. Prefixing those with a leading
# apostrophe keeps them as literal text (quotePrefix), exactly like the
# other Price cells that contain multiple dots (e.g. "38.743.25") and are
# never at risk of numeric coercion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.743.25"
$ws.Range("E2").Value = "  +1.96%  "

$ws.Range("D3").Value = "2.091.52"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'229.38"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D7").Value = "'61.08"
$ws.Range("E7").Value = "  +0.59%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +1.32%  "

$ws.Range("E10").Value = "  +0.80%  "

# Rows 11-13 shuffled rank order: WrappedliquidstakedEther2.0 moves up to
# row 11, TRON drops to row 12, Chainlink drops to row 13.
$ws.Range("B11").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C11").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D11").Value = "3.011.01"
$ws.Range("E11").Value = "  +25.65%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.104"
$ws.Range("E12").Value = "  -0.07%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'15.31"
$ws.Range("E13").Value = "  +4.74%  "

$ws.Range("D14").Value = "'22.03"
$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").Value = "'0.808"
$ws.Range("E15").Value = "  +4.67%  "

$ws.Range("D16").Value = "'5.48"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").Value = "2.100.03"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").Value = "38.713.09"
$ws.Range("E18").Value = "  +3.08%  "

$ws.Range("D19").Value = "'71.83"
$ws.Range("E19").Value = "  +2.58%  "

$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("D21").Value = "0.0₃0841"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").Value = "'227.51"
$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").Value = "'2.37"
$ws.Range("E24").Value = "  -2.40%  "

$ws.Range("D25").Value = "'2.34"
$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("D26").Value = "'171.01"
$ws.Range("E26").Value = "  +0.87%  "

$ws.Range("E27").Value = "  +1.65%  "

$ws.Range("D28").Value = "'0.139"
$ws.Range("E28").Value = "  +5.13%  "

$ws.Range("D29").Value = "'1.44"
$ws.Range("E29").Value = "  +8.87%  "

$ws.Range("D30").Value = "'19.40"
$ws.Range("E30").Value = "  +2.28%  "

$ws.Range("D31").Value = "'2.50"
$ws.Range("E31").Value = "  +5.68%  "

$ws.Range("E32").Value = "  +0.37%  "

$ws.Range("E33").Value = "  +1.88%  "

$ws.Range("D34").Value = "'4.70"
$ws.Range("E34").Value = "  +0.92%  "

$ws.Range("D35").Value = "'0.0610"
$ws.Range("E35").Value = "  +0.56%  "

$ws.Range("E36").Value = "  +1.19%  "

$ws.Range("E37").Value = "  -0.84%  "

$ws.Range("D38").Value = "'3.59"
$ws.Range("E38").Value = "  +1.84%  "

$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("D40").Value = "'17.99"
$ws.Range("E40").Value = "  +0.12%  "

$ws.Range("E41").Value = "  +4.84%  "

$ws.Range("D42").Value = "'101.05"
$ws.Range("E42").Value = "  +0.95%  "

$ws.Range("D43").Value = "1.533.26"
$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("D45").Value = "'0.0916"
$ws.Range("E45").Value = "  +1.46%  "

# Rows 46-47 swap: TrustWalletToken moves up to row 46, FraxShare drops to
# row 47.
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.13"
$ws.Range("E46").Value = "  +1.57%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'7.64"
$ws.Range("E47").Value = "  +6.14%  "

$ws.Range("D48").Value = "'4.12"
$ws.Range("E48").Value = "  -0.75%  "

$ws.Range("E49").Value = "  +1.27%  "

$ws.Range("E50").Value = "  -0.78%  "

$ws.Range("D51").Value = "2.288.51"
$ws.Range("E51").Value = "  +0.24%  "
